$d = $word.ActiveDocument

function Split-At([int]$pos) {
    # Force a run boundary at absolute character offset $pos without changing
    # any visible text, by adding then immediately deleting a zero-length
    # bookmark there (Word/this engine splits runs to host bookmark anchors,
    # and the split survives bookmark deletion).
    $bmName = "ZZsplitZZ"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($bmName, $r)
    $d.Bookmarks($bmName).Delete()
}

function Resplit-ParagraphFrom([int]$pStart, [int[]]$runLengths) {
    # Re-split the paragraph starting at $pStart back into the desired
    # run boundaries (lengths, in characters, in order).
    $offset = $pStart
    for ($i = 0; $i -lt ($runLengths.Length - 1); $i++) {
        $offset = $offset + $runLengths[$i]
        Split-At($offset)
    }
}

function Rebuild-RunGroup([int]$pStart, [int]$pEnd, [int[]]$runLengths) {
    # Rebuild the contiguous run-group spanning [$pStart, $pEnd) (which must
    # not cross a non-<w:r> sibling such as <w:proofErr/> or a bookmark) back
    # into the exact run boundaries given by $runLengths (character counts,
    # in order summing to $pEnd - $pStart). Forces an actual text edit (so
    # any stale <w:lastRenderedPageBreak/> on the group's first run is
    # dropped), then re-splits.
    $full = $d.Range($pStart, $pEnd)
    $orig = $full.Text

    $full.Text = $orig + [char]1
    $full2 = $d.Range($pStart, $pEnd + 1)
    $full2.Text = $orig

    Resplit-ParagraphFrom $pStart $runLengths
}

function Remove-LastRenderedPageBreakAndResplit-From([int]$searchFrom, [string]$findText, [int[]]$runLengths) {
    # $findText must match the *start* of the paragraph's text (used to
    # locate the paragraph via Find, searching forward from $searchFrom).
    # $runLengths are the lengths (in characters) of each desired run,
    # covering the whole paragraph text that needs to be rebuilt, in order.
    # Only safe when the whole paragraph is one contiguous run-group (no
    # <w:proofErr/> or similar siblings inside it).
    $rng = $d.Range($searchFrom, $d.Content.End)
    [void]$rng.Find.Execute($findText)
    $pStart = $rng.Start
    $para = $rng.Paragraphs(1)
    $pEnd = $para.Range.End - 1   # exclude the paragraph mark

    Rebuild-RunGroup $pStart $pEnd $runLengths
}

function Remove-LastRenderedPageBreakAndResplit([string]$findText, [int[]]$runLengths) {
    Remove-LastRenderedPageBreakAndResplit-From 0 $findText $runLengths
}

# ---------------------------------------------------------------------
# 1. Drop the stray _GoBack bookmark from the opening paragraph.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 2. Remove <w:lastRenderedPageBreak/> before "Sunday " (keep run split
#    "Sunday " | "January" | " 8, 21" | "45").
# ---------------------------------------------------------------------
Remove-LastRenderedPageBreakAndResplit "Sunday " @(7,7,6,2)

# ---------------------------------------------------------------------
# 3. Mizuno paragraph: insert "alleged " before "problem" by splitting
#    the run in three (preceding/following runs in the same paragraph
#    must be rebuilt too, since any content mutation collapses the whole
#    paragraph to a single run).
# ---------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Mizuno:  I")
$pStart = $rng.Start
$anchor = $d.Content
[void]$anchor.Find.Execute("the problem.  I’m sure we can work out a deal")
$insertPos = $anchor.Start + ("the ").Length
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertBefore("alleged ")

Resplit-ParagraphFrom $pStart @(50, 8, 137, 32, 97)

# ---------------------------------------------------------------------
# 4. Remove <w:lastRenderedPageBreak/> before "Monday " (keep run split
#    "Monday " | "January 9, 2145"). There's an earlier, unrelated
#    "Monday" earlier in the text ("...he spends Monday and Tuesday...")
#    so anchor the search after a unique nearby phrase.
# ---------------------------------------------------------------------
$anchorRng = $d.Content
[void]$anchorRng.Find.Execute("to make people’s lives better.")
Remove-LastRenderedPageBreakAndResplit-From $anchorRng.End "Monday " @(7,15)

# ---------------------------------------------------------------------
# 5a. "Matter-of-factly" paragraph: append new sentences at the very end
#     (after the <w:proofErr type="gramEnd"/> that follows "Matter-of-
#     factly."). The run-group ["Matter-of-factly","."] collapses into one
#     run when touched, so rebuild it as three runs: "Matter-of-factly",
#     ".", and the new trailing sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Matter-of-factly")
$groupStart = $rng.Start
$para = $rng.Paragraphs(1)
$mfEnd = $para.Range.End - 1   # before the paragraph mark (= right after "Matter-of-factly.", proofErr gramEnd sits here too)
$mfRange = $d.Range($mfEnd, $mfEnd)
$mfRange.InsertBefore("  He sighs a lot, and it’s kind of comical.  “I believe that the way through emotional pain is to feel it.  So that’s what I’m doing.”  There are nods all around.")

Resplit-ParagraphFrom $groupStart @(16, 1, 162)

# ---------------------------------------------------------------------
# 5b. Shelley paragraph: remove <w:lastRenderedPageBreak/> and split
#     "Shelley says: ... but way " into "Shelley " + "says: ... but way ",
#     re-inserting the _GoBack bookmark between them. This run-group ends
#     right before <w:proofErr type="gramStart"/> ("more free"), so only
#     rebuild up to there (not the whole paragraph).
# ---------------------------------------------------------------------
$rng3 = $d.Content
[void]$rng3.Find.Execute("Shelley says")
$shelleyStart = $rng3.Start
$endAnchor = $d.Content
[void]$endAnchor.Find.Execute("but way ")
$groupEnd = $endAnchor.End

Rebuild-RunGroup $shelleyStart $groupEnd @(8, 171)

$shelleySplit = $shelleyStart + 8
$d.Bookmarks.Add("_GoBack", $d.Range($shelleySplit, $shelleySplit))
